$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: TestPosts / test_get_all_posts
$ws.Cells.Item(2, 1).Value = "TestPosts"
$ws.Cells.Item(2, 2).Value = "test_get_all_posts"
$ws.Cells.Item(2, 3).Value = "Test get all posts from GET endpoint"
$ws.Cells.Item(2, 4).Value = "PASSED"
$ws.Cells.Item(2, 5).Value = 1.190232500011916
$ws.Cells.Item(2, 6).Value = "2024-05-07T13:26:54"
$ws.Cells.Item(2, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(2, 9).Value = "acceptance"

# Row 3: TestPosts / test_get_post
$ws.Cells.Item(3, 1).Value = "TestPosts"
$ws.Cells.Item(3, 2).Value = "test_get_post"
$ws.Cells.Item(3, 3).Value = "Test get a specific post object from GET endpoint"
$ws.Cells.Item(3, 4).Value = "PASSED"
$ws.Cells.Item(3, 5).Value = 0.540167900006054
$ws.Cells.Item(3, 6).Value = "2024-05-07T13:26:55"
$ws.Cells.Item(3, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(3, 9).Value = "acceptance"

# Row 4: TestPosts / test_create_post
$ws.Cells.Item(4, 1).Value = "TestPosts"
$ws.Cells.Item(4, 2).Value = "test_create_post"
$ws.Cells.Item(4, 3).Value = "Test create a new post object (posts method)"
$ws.Cells.Item(4, 4).Value = "PASSED"
$ws.Cells.Item(4, 5).Value = 1.937018500000704
$ws.Cells.Item(4, 6).Value = "2024-05-07T13:26:57"
$ws.Cells.Item(4, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(4, 9).Value = "acceptance"

# Row 5: TestPosts / test_update_post
$ws.Cells.Item(5, 1).Value = "TestPosts"
$ws.Cells.Item(5, 2).Value = "test_update_post"
$ws.Cells.Item(5, 3).Value = "Test update post object (the last created)"
$ws.Cells.Item(5, 4).Value = "PASSED"
$ws.Cells.Item(5, 5).Value = 1.737903000001097
$ws.Cells.Item(5, 6).Value = "2024-05-07T13:26:58"
$ws.Cells.Item(5, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(5, 9).Value = "acceptance"

# Row 6: TestPosts / test_delete_post
$ws.Cells.Item(6, 1).Value = "TestPosts"
$ws.Cells.Item(6, 2).Value = "test_delete_post"
$ws.Cells.Item(6, 3).Value = "Test delete a post"
$ws.Cells.Item(6, 4).Value = "PASSED"
$ws.Cells.Item(6, 5).Value = 1.653955000001588
$ws.Cells.Item(6, 6).Value = "2024-05-07T13:27:00"
$ws.Cells.Item(6, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(6, 9).Value = "acceptance"

# Row 7: TestPosts / test_required_title
$ws.Cells.Item(7, 1).Value = "TestPosts"
$ws.Cells.Item(7, 2).Value = "test_required_title"
$ws.Cells.Item(7, 3).Value = "Test title required field is not sent in request body"
$ws.Cells.Item(7, 4).Value = "PASSED"
$ws.Cells.Item(7, 5).Value = 0.5037437000137288
$ws.Cells.Item(7, 6).Value = "2024-05-07T13:27:01"
$ws.Cells.Item(7, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(7, 9).Value = "negative"

# Row 8: TestPosts / test_required_body
$ws.Cells.Item(8, 1).Value = "TestPosts"
$ws.Cells.Item(8, 2).Value = "test_required_body"
$ws.Cells.Item(8, 3).Value = "Test body required field is not sent in request body"
$ws.Cells.Item(8, 4).Value = "PASSED"
$ws.Cells.Item(8, 5).Value = 0.5080416999990121
$ws.Cells.Item(8, 6).Value = "2024-05-07T13:27:01"
$ws.Cells.Item(8, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(8, 9).Value = "negative"

# Row 9: TestPosts / test_nonexistent_post
$ws.Cells.Item(9, 1).Value = "TestPosts"
$ws.Cells.Item(9, 2).Value = "test_nonexistent_post"
$ws.Cells.Item(9, 3).Value = "Test trying to retrieve a post object that does not exist"
$ws.Cells.Item(9, 4).Value = "PASSED"
$ws.Cells.Item(9, 5).Value = 0.5132924999925308
$ws.Cells.Item(9, 6).Value = "2024-05-07T13:27:02"
$ws.Cells.Item(9, 8).Value = "gorest_api\posts\test_posts.py"
$ws.Cells.Item(9, 9).Value = "negative"

# Row 10: TestTodos / test_get_all_todos
$ws.Cells.Item(10, 1).Value = "TestTodos"
$ws.Cells.Item(10, 2).Value = "test_get_all_todos"
$ws.Cells.Item(10, 3).Value = "Test get all todos from GET endpoint"
$ws.Cells.Item(10, 4).Value = "PASSED"
$ws.Cells.Item(10, 5).Value = 1.336760799997137
$ws.Cells.Item(10, 6).Value = "2024-05-07T13:27:05"
$ws.Cells.Item(10, 8).Value = "gorest_api\todos\test_todos.py"
$ws.Cells.Item(10, 9).Value = "acceptance"

# Row 11: TestTodos / test_get_todo
$ws.Cells.Item(11, 1).Value = "TestTodos"
$ws.Cells.Item(11, 2).Value = "test_get_todo"
$ws.Cells.Item(11, 3).Value = "Test get a specific todos object from GET endpoint"
$ws.Cells.Item(11, 4).Value = "PASSED"
$ws.Cells.Item(11, 5).Value = 0.600476499996148
$ws.Cells.Item(11, 6).Value = "2024-05-07T13:27:06"
$ws.Cells.Item(11, 8).Value = "gorest_api\todos\test_todos.py"
$ws.Cells.Item(11, 9).Value = "acceptance"

# Row 12: TestTodos / test_create_todo
$ws.Cells.Item(12, 1).Value = "TestTodos"
$ws.Cells.Item(12, 2).Value = "test_create_todo"
$ws.Cells.Item(12, 3).Value = "Test create a new todo_object (posts method)"
$ws.Cells.Item(12, 4).Value = "PASSED"
$ws.Cells.Item(12, 5).Value = 1.858424400008516
$ws.Cells.Item(12, 6).Value = "2024-05-07T13:27:08"
$ws.Cells.Item(12, 8).Value = "gorest_api\todos\test_todos.py"
$ws.Cells.Item(12, 9).Value = "acceptance"

# Row 13: TestTodos / test_update_todo
$ws.Cells.Item(13, 1).Value = "TestTodos"
$ws.Cells.Item(13, 2).Value = "test_update_todo"
$ws.Cells.Item(13, 3).Value = "Test update todo_object (the last created)"
$ws.Cells.Item(13, 4).Value = "PASSED"
$ws.Cells.Item(13, 5).Value = 1.790401999998721
$ws.Cells.Item(13, 6).Value = "2024-05-07T13:27:10"
$ws.Cells.Item(13, 8).Value = "gorest_api\todos\test_todos.py"
$ws.Cells.Item(13, 9).Value = "acceptance"

# Row 14: TestTodos / test_delete_todo
$ws.Cells.Item(14, 1).Value = "TestTodos"
$ws.Cells.Item(14, 2).Value = "test_delete_todo"
$ws.Cells.Item(14, 3).Value = "Test delete a todo_object"
$ws.Cells.Item(14, 4).Value = "PASSED"
$ws.Cells.Item(14, 5).Value = 1.737185999998474
$ws.Cells.Item(14, 6).Value = "2024-05-07T13:27:11"
$ws.Cells.Item(14, 8).Value = "gorest_api\todos\test_todos.py"
$ws.Cells.Item(14, 9).Value = "acceptance"

# Row 15: TestUsers / test_get_all_users
$ws.Cells.Item(15, 1).Value = "TestUsers"
$ws.Cells.Item(15, 2).Value = "test_get_all_users"
$ws.Cells.Item(15, 3).Value = "Test get all user from GET endpoint"
$ws.Cells.Item(15, 4).Value = "PASSED"
$ws.Cells.Item(15, 5).Value = 1.26391720000538
$ws.Cells.Item(15, 6).Value = "2024-05-07T13:27:16"
$ws.Cells.Item(15, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(15, 9).Value = "acceptance"

# Row 16: TestUsers / test_get_user
$ws.Cells.Item(16, 1).Value = "TestUsers"
$ws.Cells.Item(16, 2).Value = "test_get_user"
$ws.Cells.Item(16, 3).Value = "Test get a specific user from GET endpoint"
$ws.Cells.Item(16, 4).Value = "PASSED"
$ws.Cells.Item(16, 5).Value = 0.5325499000027776
$ws.Cells.Item(16, 6).Value = "2024-05-07T13:27:16"
$ws.Cells.Item(16, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(16, 9).Value = "sanity"

# Row 17: TestUsers / test_create_user
$ws.Cells.Item(17, 1).Value = "TestUsers"
$ws.Cells.Item(17, 2).Value = "test_create_user"
$ws.Cells.Item(17, 3).Value = "Test create a new user (posts method)"
$ws.Cells.Item(17, 4).Value = "PASSED"
$ws.Cells.Item(17, 5).Value = 0.643302099997527
$ws.Cells.Item(17, 6).Value = "2024-05-07T13:27:17"
$ws.Cells.Item(17, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(17, 9).Value = "acceptance"

# Row 18: TestUsers / test_update_user
$ws.Cells.Item(18, 1).Value = "TestUsers"
$ws.Cells.Item(18, 2).Value = "test_update_user"
$ws.Cells.Item(18, 3).Value = "Test update user (the last created)"
$ws.Cells.Item(18, 4).Value = "PASSED"
$ws.Cells.Item(18, 5).Value = 1.218064999993658
$ws.Cells.Item(18, 6).Value = "2024-05-07T13:27:18"
$ws.Cells.Item(18, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(18, 9).Value = "acceptance"

# Row 19: TestUsers / test_delete_user
$ws.Cells.Item(19, 1).Value = "TestUsers"
$ws.Cells.Item(19, 2).Value = "test_delete_user"
$ws.Cells.Item(19, 3).Value = "Test delete a user"
$ws.Cells.Item(19, 4).Value = "PASSED"
$ws.Cells.Item(19, 5).Value = 1.453652500000317
$ws.Cells.Item(19, 6).Value = "2024-05-07T13:27:20"
$ws.Cells.Item(19, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(19, 9).Value = "acceptance"

# Row 20: TestUsers / test_required_field_name
$ws.Cells.Item(20, 1).Value = "TestUsers"
$ws.Cells.Item(20, 2).Value = "test_required_field_name"
$ws.Cells.Item(20, 3).Value = "Test required field is not sent in request body (name)"
$ws.Cells.Item(20, 4).Value = "PASSED"
$ws.Cells.Item(20, 5).Value = 0.5441369000036502
$ws.Cells.Item(20, 6).Value = "2024-05-07T13:27:20"
$ws.Cells.Item(20, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(20, 9).Value = "negative"

# Row 21: TestUsers / test_email_already_taken
$ws.Cells.Item(21, 1).Value = "TestUsers"
$ws.Cells.Item(21, 2).Value = "test_email_already_taken"
$ws.Cells.Item(21, 3).Value = "Test email address is already taken"
$ws.Cells.Item(21, 4).Value = "PASSED"
$ws.Cells.Item(21, 5).Value = 1.17042070000025
$ws.Cells.Item(21, 6).Value = "2024-05-07T13:27:21"
$ws.Cells.Item(21, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(21, 9).Value = "negative"

# Row 22: TestUsers / test_nonexistent_user
$ws.Cells.Item(22, 1).Value = "TestUsers"
$ws.Cells.Item(22, 2).Value = "test_nonexistent_user"
$ws.Cells.Item(22, 3).Value = "Test trying to retrieve a user that does not exist"
$ws.Cells.Item(22, 4).Value = "PASSED"
$ws.Cells.Item(22, 5).Value = 0.529216400012956
$ws.Cells.Item(22, 6).Value = "2024-05-07T13:27:22"
$ws.Cells.Item(22, 8).Value = "gorest_api\users\test_users.py"
$ws.Cells.Item(22, 9).Value = "negative"
